$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.560.06'
$ws.Range('E2').Value = '  -0.18%  '

$ws.Range('D3').Value = '1.634.86'
$ws.Range('E3').Value = '  +0.24%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.81'
$ws.Range('E5').Value = '  +0.57%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.504'
$ws.Range('E6').Value = '  +2.08%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('E8').Value = '  -0.42%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0626'
$ws.Range('E9').Value = '  +0.43%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.85'
$ws.Range('E10').Value = '  -0.48%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0842'
$ws.Range('E11').Value = '  +0.09%  '

$ws.Range('D12').Value = '1.860.56'
$ws.Range('E12').Value = '  +0.13%  '

$ws.Range('D13').Value = '1.641.35'
$ws.Range('E13').Value = '  +0.67%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.15'
$ws.Range('E14').Value = '  +1.84%  '

$ws.Range('E15').Value = '  -0.22%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.29'
$ws.Range('E16').Value = '  +3.73%  '

$ws.Range('D17').Value = '26.579.42'
$ws.Range('E17').Value = '  -0.08%  '

$ws.Range('E18').Value = '  +0.74%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.99'
$ws.Range('E19').Value = '  +3.21%  '

$ws.Range('E20').Value = '  +0.07%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.31'
$ws.Range('E21').Value = '  +0.50%  '

$ws.Range('E22').Value = '  +1.42%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.34'
$ws.Range('E23').Value = '  -0.83%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.23'
$ws.Range('E24').Value = '  +14.92%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.00'
$ws.Range('E25').Value = '  +0.03%  '

$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('E27').Value = '  -0.52%  '

$ws.Range('E28').Value = '  +0.89%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.63'
$ws.Range('E29').Value = '  +1.87%  '

$ws.Range('E30').Value = '  -1.56%  '

$ws.Range('E31').Value = '  -0.47%  '

$ws.Range('E32').Value = '  +3.61%  '

$ws.Range('E33').Value = '  +1.09%  '

$ws.Range('D34').Value = '1.266.16'
$ws.Range('E34').Value = '  +8.68%  '

$ws.Range('E35').Value = '  +0.93%  '

$ws.Range('E36').Value = '  +0.10%  '

$ws.Range('E37').Value = '  +4.07%  '

$ws.Range('E38').Value = '  +1.37%  '

$ws.Range('E39').Value = '  +0.07%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.798'
$ws.Range('E40').Value = '  -0.84%  '

$ws.Range('E41').Value = '  -2.00%  '

$ws.Range('E42').Value = '  +0.56%  '

$ws.Range('E43').Value = '  -0.38%  '

$ws.Range('D44').Value = '1.771.20'
$ws.Range('E44').Value = '  -0.01%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.32'

$ws.Range('E46').Value = '  +3.55%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.98'
$ws.Range('E47').Value = '  +0.79%  '

$ws.Range('E48').Value = '  +0.23%  '

$ws.Range('E49').Value = '  +0.34%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.61'
$ws.Range('E50').Value = '  +1.09%  '

$ws.Range('E51').Value = '  -0.61%  '
